$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header date (F1:G1) ---
$ws.Range("F1").Value2 = 41909

# --- Remove the old data rows (the "reunión #1" entry plus the 3 that follow it) and
#     the stray SUM row below the gap, then rebuild the block from scratch so no stale
#     row-level formatting (customFormat/row style) survives. ---
$ws.Rows("6:9").Delete()
$ws.Rows(7).Delete()

# Insert a blank row above the header so it lands on row 6 and row 5 is free for the
# hours-summary formula.
$ws.Rows(5).Insert()

# --- Row 5: cycle hour summary, now totalling the (moved) data rows 7:11 ---
$ws.Range("E5").Formula = "=SUM(E7:E11)/60"

# --- Row 7 (was the "Trabajé..." entry) ---
$ws.Range("A7").Value2 = 41910
$ws.Range("B7").Value2 = 0.583333333333333
$ws.Range("C7").Value2 = 0.631944444444444
$ws.Range("D7").Value2 = 25
$ws.Range("E7").Formula = "=((HOUR(C7)-HOUR(B7))*60)+(MINUTE(C7)-MINUTE(B7))-D7"
$ws.Range("F7").Value2 = 9
$ws.Range("F7").HorizontalAlignment = -4152
$ws.Range("G7").Value2 = ""
$ws.Range("H7").Value2 = "Trabajé en la creación del esquema del documento de requerimientos."

# --- Row 8 ---
$ws.Range("A8").Value2 = 41910
$ws.Range("B8").Value2 = 0.666666666666667
$ws.Range("C8").Value2 = 0.75
$ws.Range("D8").Value2 = 30
$ws.Range("E8").Formula = "=((HOUR(C8)-HOUR(B8))*60)+(MINUTE(C8)-MINUTE(B8))-D8"
$ws.Range("F8").Value2 = 12
$ws.Range("F8").HorizontalAlignment = -4152
$ws.Range("G8").Value2 = ""
$ws.Range("H8").Value2 = "Participé en el analisis de los requerimientos obtenidos en la reunión #1. Se documento una parte del diagrama de casos de uso, y falto documentar los escenarios."

# --- Row 9 ---
$ws.Range("A9").Value2 = 41911
$ws.Range("B9").Value2 = 0.340277777777778
$ws.Range("C9").Value2 = 0.395833333333333
$ws.Range("D9").Value2 = 8
$ws.Range("E9").Formula = "=((HOUR(C9)-HOUR(B9))*60)+(MINUTE(C9)-MINUTE(B9))-D9"
$ws.Range("F9").Value2 = 9
$ws.Range("F9").HorizontalAlignment = -4152
$ws.Range("G9").Value2 = ""
$ws.Range("H9").Value2 = "Terminé la creación del esquema del documento de requerimientos."

# --- Row 10 (new entry) ---
$ws.Range("A10").Value2 = 41912
$ws.Range("B10").Value2 = 0.96875
$ws.Range("C10").Value2 = 0.993055555555556
$ws.Range("D10").Value2 = 0
$ws.Range("E10").Formula = "=((HOUR(C10)-HOUR(B10))*60)+(MINUTE(C10)-MINUTE(B10))-D10"
$ws.Range("F10").HorizontalAlignment = -4152
$ws.Range("G10").Value2 = ""
$ws.Range("H10").Value2 = "Tuvimos una reunión de equipo para la asignación de las tareas del ciclo #1."

# --- Row 11 (new entry) ---
$ws.Range("A11").Value2 = 41913
$ws.Range("B11").Value2 = 0.895833333333333
$ws.Range("C11").Value2 = 0.911805555555556
$ws.Range("D11").Value2 = 0
$ws.Range("E11").Formula = "=((HOUR(C11)-HOUR(B11))*60)+(MINUTE(C11)-MINUTE(B11))-D11"
$ws.Range("F11").Value2 = 16
$ws.Range("F11").HorizontalAlignment = -4152
$ws.Range("G11").Value2 = ""
$ws.Range("H11").Value2 = "Cree la agenda para la reunión #2 con el cliente."

# --- Row heights ---
$ws.Rows(5).RowHeight = 13.75
$ws.Rows(7).RowHeight = 26.65
$ws.Rows(8).RowHeight = 39.15
$ws.Rows(9).RowHeight = 26.65
$ws.Rows(10).RowHeight = 26.65
$ws.Rows(11).RowHeight = 14.15

# --- Selection, matching the saved cursor position in the edited workbook ---
$ws.Range("F2").Select()
